# Update CDC eviction moratorium date reference from March 31, 2021 to June 30, 2021
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Until March 31, 2021, you may be eligible to stop an eviction where your landlord is asking for back rent",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Until June 30, 2021, you may be eligible to stop an eviction where your landlord is asking for back rent",
    2
)
